$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets the content previously in row 11
$ws.Range("A2").Value = 111632126
$ws.Range("B2").Value = 90687
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 5964
$ws.Range("F2").Value = "Fjällig taggsvamp s.str."
$ws.Range("G2").Value = "Sarcodon imbricatus s.str."
$ws.Range("H2").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q2").Value = 679400.8667491183
$ws.Range("R2").Value = 6612685.041705586

# Row 3 gets the content previously in row 10
$ws.Range("A3").Value = 111631628
$ws.Range("B3").Value = 89183
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 3215
$ws.Range("F3").Value = "Rödgul trumpetsvamp"
$ws.Range("G3").Value = "Craterellus lutescens"
$ws.Range("H3").Value = "(Fr.) Fr."
$ws.Range("Q3").Value = 679347.335090697
$ws.Range("R3").Value = 6612732.142182259

# Row 4 gets the content previously in row 3
$ws.Range("A4").Value = 111631615
$ws.Range("B4").Value = 90687
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = "Fjällig taggsvamp s.str."
$ws.Range("G4").Value = "Sarcodon imbricatus s.str."
$ws.Range("H4").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q4").Value = 679347.335090697
$ws.Range("R4").Value = 6612732.142182259

# Row 5 gets the content previously in row 13
$ws.Range("A5").Value = 111632936
$ws.Range("B5").Value = 90658
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 4361
$ws.Range("F5").Value = "Orange taggsvamp"
$ws.Range("G5").Value = "Hydnellum aurantiacum"
$ws.Range("H5").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q5").Value = 679389.9201578975
$ws.Range("R5").Value = 6612881.656256998

# Row 6 gets the content previously in row 4
$ws.Range("A6").Value = 111632295
$ws.Range("B6").Value = 89419
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 1204
$ws.Range("F6").Value = "Gränsticka"
$ws.Range("G6").Value = "Phellopilus nigrolimitatus"
$ws.Range("H6").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q6").Value = 679615.9172154681
$ws.Range("R6").Value = 6612689.273485693

# Row 7 gets the content previously in row 12
$ws.Range("A7").Value = 111631820
$ws.Range("B7").Value = 88819
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 5685
$ws.Range("F7").Value = "Gullgröppa"
$ws.Range("G7").Value = "Pseudomerulius aureus"
$ws.Range("H7").Value = "(Fr.) Jülich"
$ws.Range("Q7").Value = 679445.4687985049
$ws.Range("R7").Value = 6612665.387322281

# Row 10 gets the content previously in row 6
$ws.Range("A10").Value = 111631648
$ws.Range("B10").Value = 90678
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 4366
$ws.Range("F10").Value = "Skarp dropptaggsvamp"
$ws.Range("G10").Value = "Hydnellum peckii"
$ws.Range("H10").Value = "Banker"
$ws.Range("Q10").Value = 679488.9265337941
$ws.Range("R10").Value = 6612786.06067825

# Row 11 gets the content previously in row 7
$ws.Range("A11").Value = 111632616
$ws.Range("B11").Value = 89423
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = "Granticka"
$ws.Range("G11").Value = "Porodaedalea chrysoloma"
$ws.Range("H11").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q11").Value = 679697.9675722983
$ws.Range("R11").Value = 6612862.479549611

# Row 12 gets the content previously in row 5
$ws.Range("A12").Value = 111631555
$ws.Range("B12").Value = 90658
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 4361
$ws.Range("F12").Value = "Orange taggsvamp"
$ws.Range("G12").Value = "Hydnellum aurantiacum"
$ws.Range("H12").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q12").Value = 679347.335090697
$ws.Range("R12").Value = 6612732.142182259

# Row 13 gets the content previously in row 2
$ws.Range("A13").Value = 111632455
$ws.Range("B13").Value = 90687
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 5964
$ws.Range("F13").Value = "Fjällig taggsvamp s.str."
$ws.Range("G13").Value = "Sarcodon imbricatus s.str."
$ws.Range("H13").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q13").Value = 679622.1256333978
$ws.Range("R13").Value = 6612739.744341305
